$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.981.38'
$ws.Range("E2").Value = '  +1.83%  '
$ws.Range("D3").Value = '2.260.56'
$ws.Range("E3").Value = '  +1.26%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '317.65'
$ws.Range("E5").Value = '  -0.89%  '
$ws.Range("D6").Value = '101.63'
$ws.Range("E6").Value = '  +1.34%  '
$ws.Range("E7").Value = '  -1.05%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '0.555'
$ws.Range("E9").Value = '  -1.36%  '
$ws.Range("D10").Value = '37.15'
$ws.Range("E10").Value = '  +0.26%  '
$ws.Range("D11").Value = '0.0832'
$ws.Range("E11").Value = '  -0.17%  '
$ws.Range("D12").Value = '7.66'
$ws.Range("E12").Value = '  -0.28%  '
$ws.Range("E13").Value = '  -1.95%  '
$ws.Range("D14").Value = '2.608.44'
$ws.Range("E14").Value = '  +1.36%  '
$ws.Range("D15").Value = '0.860'
$ws.Range("E15").Value = '  -0.90%  '
$ws.Range("D16").Value = '14.43'
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("D17").Value = '2.264.63'
$ws.Range("E17").Value = '  +1.27%  '
$ws.Range("D18").Value = '43.907.34'
$ws.Range("E18").Value = '  +1.80%  '
$ws.Range("D19").Value = '13.51'
$ws.Range("E19").Value = '  -7.13%  '
$ws.Range("D20").Value = '0.0₃0988'
$ws.Range("E20").Value = '  +1.87%  '
$ws.Range("D21").Value = '6.56'
$ws.Range("D22").Value = '65.68'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").Value = '3.12'
$ws.Range("E23").Value = '  -1.81%  '
$ws.Range("D24").Value = '235.09'
$ws.Range("E24").Value = '  -1.49%  '
$ws.Range("D25").Value = '2.10'
$ws.Range("E25").Value = '  -3.26%  '
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("D27").Value = '10.14'
$ws.Range("E27").Value = '  +1.04%  '
$ws.Range("E28").Value = '  -3.67%  '
$ws.Range("D29").Value = '37.24'
$ws.Range("E29").Value = '  +3.46%  '
$ws.Range("D30").Value = '6.24'
$ws.Range("E30").Value = '  -2.04%  '
$ws.Range("D31").Value = '20.19'
$ws.Range("E31").Value = '  -1.07%  '
$ws.Range("D32").Value = '158.87'
$ws.Range("E32").Value = '  +3.24%  '
$ws.Range("D33").Value = '0.0853'
$ws.Range("E33").Value = '  -2.65%  '
$ws.Range("E35").Value = '  +10.35%  '
$ws.Range("E36").Value = '  +0.54%  '
$ws.Range("D37").Value = '3.07'
$ws.Range("E37").Value = '  -3.16%  '
$ws.Range("E38").Value = '  -2.29%  '
$ws.Range("D39").Value = '16.37'
$ws.Range("E39").Value = '  +20.83%  '
$ws.Range("D40").Value = '3.72'
$ws.Range("E40").Value = '  +1.49%  '
$ws.Range("D41").Value = '4.22'
$ws.Range("E41").Value = '  -5.48%  '
$ws.Range("D42").Value = '0.0316'
$ws.Range("E42").Value = '  -2.71%  '
$ws.Range("E43").Value = '  +0.26%  '
$ws.Range("D44").Value = '1.814.94'
$ws.Range("E44").Value = '  +4.46%  '
$ws.Range("D45").Value = '76.01'
$ws.Range("E45").Value = '  +0.52%  '
$ws.Range("D46").Value = '0.198'
$ws.Range("E46").Value = '  -3.46%  '
$ws.Range("D47").Value = '82.40'
$ws.Range("E47").Value = '  -4.19%  '
$ws.Range("E48").Value = '  -1.86%  '
$ws.Range("D49").Value = '105.02'
$ws.Range("E49").Value = '  +1.86%  '
$ws.Range("D50").Value = '58.40'
$ws.Range("E50").Value = '  -0.82%  '
$ws.Range("D51").Value = '1.67'
$ws.Range("E51").Value = '  +5.53%  '
